$wb = $excel.ActiveWorkbook

# --- Sheet "Calc" updates (row 3) ---
$calc = $wb.Worksheets.Item("Calc")
$calc.Range("AP3").Value = 0.5389
$calc.Range("AQ3").Value = 0.190826111697545
$calc.Range("AW3").Value = 0.5361325914
$calc.Range("AY3").Value = 0.1898816654737066
$calc.Range("BC3").Value = 0.5689426994882338

# --- Sheet "Results" updates (row 3) ---
$results = $wb.Worksheets.Item("Results")
$results.Range("N3").Value = 0.5389
$results.Range("P3").Value = 0.5361325914
$results.Range("R3").Value = 0.5689426994882338
